$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
  2  = @(0.1169995834814548, 0.3048912486333797, 0.1496068669990043, 0.5333859586016987, 1.104883657715537)
  3  = @(0.1169995834814548, 1.626987699542094,  0.1496068669990043, 0.5333859586016987, 2.426980108624251)
  4  = @(1.445647641019636,  1.626987699542094,  0.7210945179870265, 0.5333859586016987, 4.327115817150455)
  5  = @(3.272327238179451,  1.626987699542094,  0.1496068669990043, 0.5333859586016987, 5.582307763322248)
  6  = @(0.003078177322033415, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 2.31305870246483)
  7  = @(3.272327238179451,  1.626987699542094,  0.7210945179870265, 0.5333859586016987, 6.15379541431027)
  8  = @(0.04172184405617529, 0.04103571897497393, 0.1496068669990043, 0.5333859586016987, 0.7657503886318522)
  9  = @(1.445647641019636,  1.626987699542094,  18.71679738969934,  0.5333859586016987, 22.32281868886277)
  10 = @(1.445647641019636,  1.626987699542094,  3.223369029078222, 0.5333859586016987, 6.82939032824165)
  11 = @(3.272327238179451,  1.626987699542094,  3.223369029078222, 0.5333859586016987, 8.656069925401464)
  12 = @(3.272327238179451,  1.626987699542094,  0.1496068669990043, 0.5333859586016987, 5.582307763322248)
  13 = @(3.272327238179451,  1.626987699542094,  0.7210945179870265, 0.5333859586016987, 6.15379541431027)
  14 = @(0.6545652718822623, 1.626987699542094,  0.7210945179870265, 0.5333859586016987, 3.536033448013082)
  15 = @(3.272327238179451,  1.626987699542094,  0.7210945179870265, 0.5333859586016987, 6.15379541431027)
  16 = @(3.272327238179451,  1.626987699542094,  0.1496068669990043, 0.5333859586016987, 5.582307763322248)
  17 = @(0.01253208636536152, 0.04103571897497393, 0.7210945179870265, 0.5333859586016987, 1.308048281929061)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 7).Value = $vals[4]
}
